$d = $word.ActiveDocument

# Locate the paragraph that ends with "LOB1012: Estatística (Requisito fraco)".
# This paragraph must be kept intact; we only remove what follows it.
$lob = $d.Content
$found = $lob.Find.Execute("LOB1012: Estatística (Requisito fraco)", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)

# Locate the copyright/footer paragraph that must be removed, along with the
# blank paragraph and the page-break paragraph sitting between it and the
# "LOB1012" paragraph above.
$copyright = $d.Content
$found2 = $copyright.Find.Execute(
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Expand both hits by one character so they include their trailing paragraph
# mark, then delete everything from just after the "LOB1012" paragraph mark
# through the end (incl. paragraph mark) of the copyright paragraph. That
# removes the blank paragraph, the page-break paragraph, and the copyright
# paragraph in one shot while leaving the "LOB1012" paragraph (and the two
# paragraphs that originally followed the copyright notice) untouched.
$deleteStart = $lob.End + 1
$deleteEnd = $copyright.End + 1

$toRemove = $d.Range($deleteStart, $deleteEnd)
$toRemove.Delete()
